# Update the manned ship test output sheet with the summary-of-readings
# results: columns B, C, D, E and F get new truth-state values for rows
# 6 through 86 (time steps 4 through 84).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (System 1 Truth State): rows 6-24 -> 0
$ws.Range("B6:B24").Value = 0

# Column C (System 2 Truth State): rows 6-80 -> 2
$ws.Range("C6:C80").Value = 2

# Column D (System 3 Truth State): piecewise updates
$ws.Range("D18:D25").Value = 2
$ws.Range("D27:D59").Value = 1
$ws.Range("D60:D71").Value = 2
$ws.Range("D72:D72").Value = 1
$ws.Range("D73:D78").Value = 2
$ws.Range("D79:D86").Value = 1

# Column E (System 4 Truth State): rows 6-19 -> 0
$ws.Range("E6:E19").Value = 0

# Column F (Ship Truth State): rows 6-24 -> 0
$ws.Range("F6:F24").Value = 0
